$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @(
    @{ Addr = 'D2'; Value = '56.786.95' }
    @{ Addr = 'E2'; Value = '  +1.59%  ' }
    @{ Addr = 'D3'; Value = '2.501.89' }
    @{ Addr = 'E3'; Value = '  -0.51%  ' }
    @{ Addr = 'E4'; Value = '  +0.01%  ' }
    @{ Addr = 'D5'; Value = '494.89' }
    @{ Addr = 'E5'; Value = '  +0.73%  ' }
    @{ Addr = 'D6'; Value = '153.19' }
    @{ Addr = 'E6'; Value = '  +8.61%  ' }
    @{ Addr = 'D7'; Value = '0.996' }
    @{ Addr = 'E7'; Value = '  -0.37%  ' }
    @{ Addr = 'D8'; Value = '0.515' }
    @{ Addr = 'E8'; Value = '  +0.39%  ' }
    @{ Addr = 'D9'; Value = '2.513.01' }
    @{ Addr = 'E9'; Value = '  +0.08%  ' }
    @{ Addr = 'D10'; Value = '5.75' }
    @{ Addr = 'E10'; Value = '  +3.23%  ' }
    @{ Addr = 'D11'; Value = '0.0994' }
    @{ Addr = 'E11'; Value = '  +0.33%  ' }
    @{ Addr = 'D12'; Value = '0.336' }
    @{ Addr = 'E12'; Value = '  +0.88%  ' }
    @{ Addr = 'D14'; Value = '2.937.63' }
    @{ Addr = 'E14'; Value = '  -0.49%  ' }
    @{ Addr = 'D15'; Value = '56.812.27' }
    @{ Addr = 'E15'; Value = '  +1.68%  ' }
    @{ Addr = 'D16'; Value = '21.35' }
    @{ Addr = 'E16'; Value = '  +1.89%  ' }
    @{ Addr = 'D17'; Value = '0.0000138' }
    @{ Addr = 'E17'; Value = '  -0.47%  ' }
    @{ Addr = 'D18'; Value = '2.506.36' }
    @{ Addr = 'E18'; Value = '  -0.55%  ' }
    @{ Addr = 'D19'; Value = '4.58' }
    @{ Addr = 'E19'; Value = '  +3.41%  ' }
    @{ Addr = 'D20'; Value = '10.34' }
    @{ Addr = 'E20'; Value = '  +2.54%  ' }
    @{ Addr = 'D21'; Value = '322.67' }
    @{ Addr = 'E21'; Value = '  -0.08%  ' }
    @{ Addr = 'E22'; Value = '  -0.18%  ' }
    @{ Addr = 'E23'; Value = '  +1.82%  ' }
    @{ Addr = 'D24'; Value = '58.66' }
    @{ Addr = 'E24'; Value = '  +0.38%  ' }
    @{ Addr = 'D25'; Value = '0.411' }
    @{ Addr = 'E25'; Value = '  -0.50%  ' }
    @{ Addr = 'D26'; Value = '1.00' }
    @{ Addr = 'E26'; Value = '  -0.72%  ' }
    @{ Addr = 'E27'; Value = '  -5.60%  ' }
    @{ Addr = 'D28'; Value = '2.602.75' }
    @{ Addr = 'E28'; Value = '  -0.83%  ' }
    @{ Addr = 'D29'; Value = '7.60' }
    @{ Addr = 'E29'; Value = '  +1.80%  ' }
    @{ Addr = 'D30'; Value = '0.0₃0811' }
    @{ Addr = 'E30'; Value = '  +1.48%  ' }
    @{ Addr = 'D31'; Value = '0.999' }
    @{ Addr = 'E31'; Value = '  -0.08%  ' }
    @{ Addr = 'D32'; Value = '151.46' }
    @{ Addr = 'E32'; Value = '  +0.97%  ' }
    @{ Addr = 'D33'; Value = '18.35' }
    @{ Addr = 'E33'; Value = '  +0.24%  ' }
    @{ Addr = 'D34'; Value = '1.53' }
    @{ Addr = 'E34'; Value = '  +1.34%  ' }
    @{ Addr = 'D35'; Value = '5.30' }
    @{ Addr = 'E35'; Value = '  +1.18%  ' }
    @{ Addr = 'E36'; Value = '  +2.44%  ' }
    @{ Addr = 'D37'; Value = '3.81' }
    @{ Addr = 'E37'; Value = '  +1.89%  ' }
    @{ Addr = 'D38'; Value = '0.876' }
    @{ Addr = 'E38'; Value = '  -0.95%  ' }
    @{ Addr = 'E39'; Value = '  +5.09%  ' }
    @{ Addr = 'D40'; Value = '34.22' }
    @{ Addr = 'E40'; Value = '  -0.60%  ' }
    @{ Addr = 'E41'; Value = '  +2.14%  ' }
    @{ Addr = 'E42'; Value = '  +1.39%  ' }
    @{ Addr = 'E43'; Value = '  +0.40%  ' }
    @{ Addr = 'D44'; Value = '0.994' }
    @{ Addr = 'E44'; Value = '  -0.44%  ' }
    @{ Addr = 'B45'; Value = 'Bittensor' }
    @{ Addr = 'C45'; Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao' }
    @{ Addr = 'D45'; Value = '270.98' }
    @{ Addr = 'E45'; Value = '  +4.15%  ' }
    @{ Addr = 'B46'; Value = 'RenderToken' }
    @{ Addr = 'C46'; Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr' }
    @{ Addr = 'D46'; Value = '4.90' }
    @{ Addr = 'E46'; Value = '  +2.36%  ' }
    @{ Addr = 'D47'; Value = '0.0937' }
    @{ Addr = 'E47'; Value = '  +2.62%  ' }
    @{ Addr = 'D48'; Value = '0.0230' }
    @{ Addr = 'E48'; Value = '  +1.42%  ' }
    @{ Addr = 'D49'; Value = '10.20' }
    @{ Addr = 'E49'; Value = '  +0.51%  ' }
    @{ Addr = 'D50'; Value = '18.01' }
    @{ Addr = 'E50'; Value = '  +2.30%  ' }
    @{ Addr = 'D51'; Value = '1.905.60' }
    @{ Addr = 'E51'; Value = '  -4.91%  ' }
)

foreach ($item in $changes) {
    $cell = $ws.Range($item.Addr)
    $cell.NumberFormat = "@"
    $cell.Value = $item.Value
}
